# Fix ambiguous formulation: "load balance" -> "load inbalance" (sic, per
# the commit's wording), splitting the run the way PowerPoint's own
# proofing pass would once it has retyped a word it does not recognize.
#
# Slide 71, "Content Placeholder 2" shape, the paragraph (outline level 2,
# i.e. lvl="1") that currently reads:
#   "Unexpected load balance between processes/threads"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(71)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$para = $tr.Paragraphs(8, 1)

# Sanity check - make sure we are editing the expected paragraph before
# touching anything. (Paragraphs(...).Text carries a trailing CR for the
# paragraph mark, so trim it before comparing.)
$paraText = $para.Text.TrimEnd("`r")
if ($paraText -eq "Unexpected load balance between processes/threads") {

    # Replace just the word "balance" -> "inbalance". Re-assigning a
    # sub-range's .Text creates a fresh run boundary around the replaced
    # text, which is exactly the split the diff shows.
    $word = $para.Characters(17, 7)
    $word.Text = "inbalance"

    # Re-assert the run boundaries for the trailing " between
    # processes/threads" so the single space also becomes its own run,
    # matching the four-run layout in the target slide.
    $space = $para.Characters(26, 1)
    $space.Text = " "

    $rest = $para.Characters(27, 25)
    $rest.Text = "between processes/threads"
}
